$wb = $excel.ActiveWorkbook

$wsSched = $wb.Worksheets.Item("Schedule")

# --- Update existing row 2 ---
$wsSched.Cells.Item(2, 5).Value = 335.6593695
$wsSched.Cells.Item(2, 6).Value = 8.072615909090908

# --- New row 3 ---
$wsSched.Cells.Item(3, 1).Value = 46071.0625
$wsSched.Cells.Item(3, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsSched.Cells.Item(3, 2).Value = 46071.22916666666
$wsSched.Cells.Item(3, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsSched.Cells.Item(3, 3).Value = 4
$wsSched.Cells.Item(3, 4).Value = 15.12
$wsSched.Cells.Item(3, 5).Value = 544.72509975
$wsSched.Cells.Item(3, 6).Value = 36.02679231150794

# --- New row 4 ---
$wsSched.Cells.Item(4, 1).Value = 46071.3125
$wsSched.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsSched.Cells.Item(4, 2).Value = 46071.60416666666
$wsSched.Cells.Item(4, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsSched.Cells.Item(4, 3).Value = 7
$wsSched.Cells.Item(4, 4).Value = 26.46
$wsSched.Cells.Item(4, 5).Value = 94.920891
$wsSched.Cells.Item(4, 6).Value = 3.587335260770975

$ws = $wb.Worksheets.Item("Detailed")

# --- Updates to existing rows (14-49) ---
$ws.Cells.Item(14, 2).Value = 76.36685
$ws.Cells.Item(15, 2).Value = 73.45554
$ws.Cells.Item(15, 3).Value = "historical"
$ws.Cells.Item(16, 2).Value = 56.98
$ws.Cells.Item(16, 3).Value = "historical"
$ws.Cells.Item(17, 2).Value = 28.67198
$ws.Cells.Item(17, 3).Value = "historical"
$ws.Cells.Item(18, 2).Value = 14.66346
$ws.Cells.Item(18, 3).Value = "historical"
$ws.Cells.Item(19, 2).Value = 0.65347
$ws.Cells.Item(19, 3).Value = "historical"
$ws.Cells.Item(20, 2).Value = 0.07736
$ws.Cells.Item(20, 3).Value = "historical"
$ws.Cells.Item(21, 2).Value = -4.28635
$ws.Cells.Item(21, 3).Value = "historical"
$ws.Cells.Item(22, 2).Value = -5.25609
$ws.Cells.Item(22, 3).Value = "historical"
$ws.Cells.Item(23, 2).Value = -5.97948
$ws.Cells.Item(23, 3).Value = "historical"
$ws.Cells.Item(24, 2).Value = -8.474729999999999
$ws.Cells.Item(24, 3).Value = "historical"
$ws.Cells.Item(25, 2).Value = -8.31752
$ws.Cells.Item(25, 3).Value = "historical"
$ws.Cells.Item(26, 2).Value = -6.24252
$ws.Cells.Item(26, 3).Value = "historical"
$ws.Cells.Item(27, 2).Value = -7.97915
$ws.Cells.Item(27, 3).Value = "historical"
$ws.Cells.Item(28, 2).Value = -7.73511
$ws.Cells.Item(28, 3).Value = "historical"
$ws.Cells.Item(29, 2).Value = -7.2788
$ws.Cells.Item(29, 3).Value = "historical"
$ws.Cells.Item(30, 3).Value = "historical"
$ws.Cells.Item(31, 2).Value = -2.17044
$ws.Cells.Item(31, 3).Value = "historical"
$ws.Cells.Item(32, 2).Value = -4
$ws.Cells.Item(32, 3).Value = "historical"
$ws.Cells.Item(33, 2).Value = 27.08092
$ws.Cells.Item(33, 3).Value = "historical"
$ws.Cells.Item(34, 2).Value = 54.31663
$ws.Cells.Item(35, 2).Value = 56.29041
$ws.Cells.Item(36, 2).Value = 61.63378
$ws.Cells.Item(37, 2).Value = 84.79000000000001
$ws.Cells.Item(39, 2).Value = 68.79849
$ws.Cells.Item(40, 2).Value = 81.76921
$ws.Cells.Item(41, 2).Value = 102.70645
$ws.Cells.Item(42, 2).Value = 87.70733
$ws.Cells.Item(43, 2).Value = 82.8614
$ws.Cells.Item(44, 2).Value = 95.05185
$ws.Cells.Item(45, 2).Value = 79.98025
$ws.Cells.Item(49, 2).Value = 77.0551

# --- New rows (50-97) ---
$ws.Cells.Item(50, 1).Value = 46071
$ws.Cells.Item(50, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 2).Value = 78.69226999999999
$ws.Cells.Item(50, 3).Value = "forecast"
$ws.Cells.Item(50, 4).Value = 46071
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(50, 5).Value = "OFF"

$ws.Cells.Item(51, 1).Value = 46071.02083333334
$ws.Cells.Item(51, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(51, 2).Value = 79.95
$ws.Cells.Item(51, 3).Value = "forecast"
$ws.Cells.Item(51, 4).Value = 46071
$ws.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(51, 5).Value = "OFF"

$ws.Cells.Item(52, 1).Value = 46071.04166666666
$ws.Cells.Item(52, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(52, 2).Value = 79.13376
$ws.Cells.Item(52, 3).Value = "forecast"
$ws.Cells.Item(52, 4).Value = 46071
$ws.Cells.Item(52, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(52, 5).Value = "OFF"

$ws.Cells.Item(53, 1).Value = 46071.0625
$ws.Cells.Item(53, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(53, 2).Value = 79.20332999999999
$ws.Cells.Item(53, 3).Value = "forecast"
$ws.Cells.Item(53, 4).Value = 46071
$ws.Cells.Item(53, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(53, 5).Value = "ON"

$ws.Cells.Item(54, 1).Value = 46071.08333333334
$ws.Cells.Item(54, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 2).Value = 77.70349
$ws.Cells.Item(54, 3).Value = "forecast"
$ws.Cells.Item(54, 4).Value = 46071
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(54, 5).Value = "ON"

$ws.Cells.Item(55, 1).Value = 46071.10416666666
$ws.Cells.Item(55, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 2).Value = 65.81824
$ws.Cells.Item(55, 3).Value = "forecast"
$ws.Cells.Item(55, 4).Value = 46071
$ws.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(55, 5).Value = "ON"

$ws.Cells.Item(56, 1).Value = 46071.125
$ws.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 2).Value = 65.45835
$ws.Cells.Item(56, 3).Value = "forecast"
$ws.Cells.Item(56, 4).Value = 46071
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(56, 5).Value = "ON"

$ws.Cells.Item(57, 1).Value = 46071.14583333334
$ws.Cells.Item(57, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(57, 2).Value = 65.06932
$ws.Cells.Item(57, 3).Value = "forecast"
$ws.Cells.Item(57, 4).Value = 46071
$ws.Cells.Item(57, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(57, 5).Value = "ON"

$ws.Cells.Item(58, 1).Value = 46071.16666666666
$ws.Cells.Item(58, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(58, 2).Value = 65.05297
$ws.Cells.Item(58, 3).Value = "forecast"
$ws.Cells.Item(58, 4).Value = 46071
$ws.Cells.Item(58, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(58, 5).Value = "ON"

$ws.Cells.Item(59, 1).Value = 46071.1875
$ws.Cells.Item(59, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 2).Value = 65
$ws.Cells.Item(59, 3).Value = "forecast"
$ws.Cells.Item(59, 4).Value = 46071
$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(59, 5).Value = "ON"

$ws.Cells.Item(60, 1).Value = 46071.20833333334
$ws.Cells.Item(60, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 2).Value = 75.38670999999999
$ws.Cells.Item(60, 3).Value = "forecast"
$ws.Cells.Item(60, 4).Value = 46071
$ws.Cells.Item(60, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(60, 5).Value = "ON"

$ws.Cells.Item(61, 1).Value = 46071.22916666666
$ws.Cells.Item(61, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(61, 2).Value = 79.9502
$ws.Cells.Item(61, 3).Value = "forecast"
$ws.Cells.Item(61, 4).Value = 46071
$ws.Cells.Item(61, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(61, 5).Value = "OFF"

$ws.Cells.Item(62, 1).Value = 46071.25
$ws.Cells.Item(62, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 2).Value = 105.79
$ws.Cells.Item(62, 3).Value = "forecast"
$ws.Cells.Item(62, 4).Value = 46071
$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(62, 5).Value = "OFF"

$ws.Cells.Item(63, 1).Value = 46071.27083333334
$ws.Cells.Item(63, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 2).Value = 106.59569
$ws.Cells.Item(63, 3).Value = "forecast"
$ws.Cells.Item(63, 4).Value = 46071
$ws.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(63, 5).Value = "OFF"

$ws.Cells.Item(64, 1).Value = 46071.29166666666
$ws.Cells.Item(64, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 2).Value = 84.79000000000001
$ws.Cells.Item(64, 3).Value = "forecast"
$ws.Cells.Item(64, 4).Value = 46071
$ws.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(64, 5).Value = "OFF"

$ws.Cells.Item(65, 1).Value = 46071.3125
$ws.Cells.Item(65, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(65, 2).Value = 29.37442
$ws.Cells.Item(65, 3).Value = "forecast"
$ws.Cells.Item(65, 4).Value = 46071
$ws.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(65, 5).Value = "ON"

$ws.Cells.Item(66, 1).Value = 46071.33333333334
$ws.Cells.Item(66, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(66, 2).Value = 3.12401
$ws.Cells.Item(66, 3).Value = "forecast"
$ws.Cells.Item(66, 4).Value = 46071
$ws.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(66, 5).Value = "ON"

$ws.Cells.Item(67, 1).Value = 46071.35416666666
$ws.Cells.Item(67, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(67, 2).Value = 0.50973
$ws.Cells.Item(67, 3).Value = "forecast"
$ws.Cells.Item(67, 4).Value = 46071
$ws.Cells.Item(67, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(67, 5).Value = "ON"

$ws.Cells.Item(68, 1).Value = 46071.375
$ws.Cells.Item(68, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(68, 2).Value = 0.51
$ws.Cells.Item(68, 3).Value = "forecast"
$ws.Cells.Item(68, 4).Value = 46071
$ws.Cells.Item(68, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(68, 5).Value = "ON"

$ws.Cells.Item(69, 1).Value = 46071.39583333334
$ws.Cells.Item(69, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(69, 2).Value = -2.28597
$ws.Cells.Item(69, 3).Value = "forecast"
$ws.Cells.Item(69, 4).Value = 46071
$ws.Cells.Item(69, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(69, 5).Value = "ON"

$ws.Cells.Item(70, 1).Value = 46071.41666666666
$ws.Cells.Item(70, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(70, 2).Value = 0.38444
$ws.Cells.Item(70, 3).Value = "forecast"
$ws.Cells.Item(70, 4).Value = 46071
$ws.Cells.Item(70, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(70, 5).Value = "ON"

$ws.Cells.Item(71, 1).Value = 46071.4375
$ws.Cells.Item(71, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(71, 2).Value = 0.51
$ws.Cells.Item(71, 3).Value = "forecast"
$ws.Cells.Item(71, 4).Value = 46071
$ws.Cells.Item(71, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(71, 5).Value = "ON"

$ws.Cells.Item(72, 1).Value = 46071.45833333334
$ws.Cells.Item(72, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(72, 2).Value = -2.65692
$ws.Cells.Item(72, 3).Value = "forecast"
$ws.Cells.Item(72, 4).Value = 46071
$ws.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(72, 5).Value = "ON"

$ws.Cells.Item(73, 1).Value = 46071.47916666666
$ws.Cells.Item(73, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(73, 2).Value = 0.51
$ws.Cells.Item(73, 3).Value = "forecast"
$ws.Cells.Item(73, 4).Value = 46071
$ws.Cells.Item(73, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(73, 5).Value = "ON"

$ws.Cells.Item(74, 1).Value = 46071.5
$ws.Cells.Item(74, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 2).Value = 2.80179
$ws.Cells.Item(74, 3).Value = "forecast"
$ws.Cells.Item(74, 4).Value = 46071
$ws.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(74, 5).Value = "ON"

$ws.Cells.Item(75, 1).Value = 46071.52083333334
$ws.Cells.Item(75, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 2).Value = 3.4829
$ws.Cells.Item(75, 3).Value = "forecast"
$ws.Cells.Item(75, 4).Value = 46071
$ws.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(75, 5).Value = "ON"

$ws.Cells.Item(76, 1).Value = 46071.54166666666
$ws.Cells.Item(76, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(76, 2).Value = 3.83553
$ws.Cells.Item(76, 3).Value = "forecast"
$ws.Cells.Item(76, 4).Value = 46071
$ws.Cells.Item(76, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(76, 5).Value = "ON"

$ws.Cells.Item(77, 1).Value = 46071.5625
$ws.Cells.Item(77, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77, 2).Value = 21.19483
$ws.Cells.Item(77, 3).Value = "forecast"
$ws.Cells.Item(77, 4).Value = 46071
$ws.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(77, 5).Value = "ON"

$ws.Cells.Item(78, 1).Value = 46071.58333333334
$ws.Cells.Item(78, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78, 2).Value = 36.06
$ws.Cells.Item(78, 3).Value = "forecast"
$ws.Cells.Item(78, 4).Value = 46071
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(78, 5).Value = "ON"

$ws.Cells.Item(79, 1).Value = 46071.60416666666
$ws.Cells.Item(79, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79, 2).Value = 57.12018
$ws.Cells.Item(79, 3).Value = "forecast"
$ws.Cells.Item(79, 4).Value = 46071
$ws.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(79, 5).Value = "OFF"

$ws.Cells.Item(80, 1).Value = 46071.625
$ws.Cells.Item(80, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(80, 2).Value = 78
$ws.Cells.Item(80, 3).Value = "forecast"
$ws.Cells.Item(80, 4).Value = 46071
$ws.Cells.Item(80, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(80, 5).Value = "OFF"

$ws.Cells.Item(81, 1).Value = 46071.64583333334
$ws.Cells.Item(81, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81, 2).Value = 71.85057999999999
$ws.Cells.Item(81, 3).Value = "forecast"
$ws.Cells.Item(81, 4).Value = 46071
$ws.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(81, 5).Value = "OFF"

$ws.Cells.Item(82, 1).Value = 46071.66666666666
$ws.Cells.Item(82, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(82, 2).Value = 69.80280999999999
$ws.Cells.Item(82, 3).Value = "forecast"
$ws.Cells.Item(82, 4).Value = 46071
$ws.Cells.Item(82, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(82, 5).Value = "OFF"

$ws.Cells.Item(83, 1).Value = 46071.6875
$ws.Cells.Item(83, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(83, 2).Value = 57.06008
$ws.Cells.Item(83, 3).Value = "forecast"
$ws.Cells.Item(83, 4).Value = 46071
$ws.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(83, 5).Value = "OFF"

$ws.Cells.Item(84, 1).Value = 46071.70833333334
$ws.Cells.Item(84, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84, 2).Value = 85.40304
$ws.Cells.Item(84, 3).Value = "forecast"
$ws.Cells.Item(84, 4).Value = 46071
$ws.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(84, 5).Value = "OFF"

$ws.Cells.Item(85, 1).Value = 46071.72916666666
$ws.Cells.Item(85, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(85, 2).Value = 57.03889
$ws.Cells.Item(85, 3).Value = "forecast"
$ws.Cells.Item(85, 4).Value = 46071
$ws.Cells.Item(85, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(85, 5).Value = "OFF"

$ws.Cells.Item(86, 1).Value = 46071.75
$ws.Cells.Item(86, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 2).Value = 67.55777999999999
$ws.Cells.Item(86, 3).Value = "forecast"
$ws.Cells.Item(86, 4).Value = 46071
$ws.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(86, 5).Value = "OFF"

$ws.Cells.Item(87, 1).Value = 46071.77083333334
$ws.Cells.Item(87, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(87, 2).Value = 93.54431
$ws.Cells.Item(87, 3).Value = "forecast"
$ws.Cells.Item(87, 4).Value = 46071
$ws.Cells.Item(87, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(87, 5).Value = "OFF"

$ws.Cells.Item(88, 1).Value = 46071.79166666666
$ws.Cells.Item(88, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(88, 2).Value = 150.62914
$ws.Cells.Item(88, 3).Value = "forecast"
$ws.Cells.Item(88, 4).Value = 46071
$ws.Cells.Item(88, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(88, 5).Value = "OFF"

$ws.Cells.Item(89, 1).Value = 46071.8125
$ws.Cells.Item(89, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(89, 2).Value = 167.69046
$ws.Cells.Item(89, 3).Value = "forecast"
$ws.Cells.Item(89, 4).Value = 46071
$ws.Cells.Item(89, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(89, 5).Value = "OFF"

$ws.Cells.Item(90, 1).Value = 46071.83333333334
$ws.Cells.Item(90, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(90, 2).Value = 142.69675
$ws.Cells.Item(90, 3).Value = "forecast"
$ws.Cells.Item(90, 4).Value = 46071
$ws.Cells.Item(90, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(90, 5).Value = "OFF"

$ws.Cells.Item(91, 1).Value = 46071.85416666666
$ws.Cells.Item(91, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(91, 2).Value = 129.77149
$ws.Cells.Item(91, 3).Value = "forecast"
$ws.Cells.Item(91, 4).Value = 46071
$ws.Cells.Item(91, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(91, 5).Value = "OFF"

$ws.Cells.Item(92, 1).Value = 46071.875
$ws.Cells.Item(92, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 2).Value = 154.2
$ws.Cells.Item(92, 3).Value = "forecast"
$ws.Cells.Item(92, 4).Value = 46071
$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(92, 5).Value = "OFF"

$ws.Cells.Item(93, 1).Value = 46071.89583333334
$ws.Cells.Item(93, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(93, 2).Value = 111.84297
$ws.Cells.Item(93, 3).Value = "forecast"
$ws.Cells.Item(93, 4).Value = 46071
$ws.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(93, 5).Value = "OFF"

$ws.Cells.Item(94, 1).Value = 46071.91666666666
$ws.Cells.Item(94, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(94, 2).Value = 105.9578
$ws.Cells.Item(94, 3).Value = "forecast"
$ws.Cells.Item(94, 4).Value = 46071
$ws.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(94, 5).Value = "OFF"

$ws.Cells.Item(95, 1).Value = 46071.9375
$ws.Cells.Item(95, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(95, 2).Value = 105.37785
$ws.Cells.Item(95, 3).Value = "forecast"
$ws.Cells.Item(95, 4).Value = 46071
$ws.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(95, 5).Value = "OFF"

$ws.Cells.Item(96, 1).Value = 46071.95833333334
$ws.Cells.Item(96, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(96, 2).Value = 97.21102999999999
$ws.Cells.Item(96, 3).Value = "forecast"
$ws.Cells.Item(96, 4).Value = 46071
$ws.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(96, 5).Value = "OFF"

$ws.Cells.Item(97, 1).Value = 46071.97916666666
$ws.Cells.Item(97, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(97, 2).Value = 89.82371000000001
$ws.Cells.Item(97, 3).Value = "forecast"
$ws.Cells.Item(97, 4).Value = 46071
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(97, 5).Value = "OFF"
